$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.0498220640569395
$ws1.Range("C2").Value = 0.0498220640569395
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.09491525423728814
$ws1.Range("F2").Value = 0.2077151335311573
$ws1.Range("G2").Value = 0.5768621236133122
$ws1.Range("H2").Value = 0.6927166934189406
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 534
$ws1.Range("K2").Value = 0
$ws1.Range("L2").Value = 0

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0

$ws2.Range("B3").Value = 0.0498220640569395
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.09491525423728814

$ws2.Range("B4").Value = 0.0498220640569395
$ws2.Range("C4").Value = 0.0498220640569395
$ws2.Range("D4").Value = 0.0498220640569395
$ws2.Range("E4").Value = 0.0498220640569395

$ws2.Range("B5").Value = 0.02491103202846975
$ws2.Range("C5").Value = 0.5
$ws2.Range("D5").Value = 0.04745762711864407

$ws2.Range("B6").Value = 0.002482238066893783
$ws2.Range("C6").Value = 0.0498220640569395
$ws2.Range("D6").Value = 0.004728873876590867

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 534

$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
